# Auto-generated edit script applying scheduled market-data refresh to the Asura_Profits workbook.
# For each affected leve row, updates price/profit columns (H-N) to the latest scrape values.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""
$ws.Range("H129").Value = 916.6
$ws.Range("J129").Value = 1048.2609
$ws.Range("L129").Value = 3144.7827
$ws.Range("N129").Value = -13144.7827
$ws.Range("H132").Value = 2964.4546
$ws.Range("I132").Value = 2577.926
$ws.Range("J132").Value = 3578.353
$ws.Range("K132").Value = 7733.778
$ws.Range("L132").Value = 10735.059
$ws.Range("M132").Value = -5203.778
$ws.Range("N132").Value = -15795.059
$ws.Range("H137").Value = 1280.75
$ws.Range("I137").Value = 1039.7188
$ws.Range("J137").Value = 1762.8125
$ws.Range("K137").Value = 3119.1564
$ws.Range("L137").Value = 5288.4375
$ws.Range("M137").Value = -569.1564000000003
$ws.Range("N137").Value = -10388.4375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1089.9354
$ws.Range("I74").Value = 869
$ws.Range("J74").Value = 1553.9
$ws.Range("K74").Value = 869
$ws.Range("L74").Value = 1553.9
$ws.Range("M74").Value = 5
$ws.Range("N74").Value = -3301.9
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H77").Value = 1089.9354
$ws.Range("I77").Value = 869
$ws.Range("J77").Value = 1553.9
$ws.Range("K77").Value = 4345
$ws.Range("L77").Value = 7769.5
$ws.Range("M77").Value = 23
$ws.Range("N77").Value = -16505.5
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 7507.5
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 10015
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 10015
$ws.Range("M21").Value = -4765
$ws.Range("N21").Value = -10485
$ws.Range("H41").Value = 17091.25
$ws.Range("I41").Value = 1650
$ws.Range("K41").Value = 1650
$ws.Range("M41").Value = -1222
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = ""
$ws.Range("N50").Value = ""
$ws.Range("H51").Value = 29861.75
$ws.Range("J51").Value = 29861.75
$ws.Range("L51").Value = 29861.75
$ws.Range("N51").Value = -31333.75
$ws.Range("H60").Value = 36468.668
$ws.Range("I60").Value = 20000
$ws.Range("J60").Value = 39762.4
$ws.Range("K60").Value = 20000
$ws.Range("L60").Value = 39762.4
$ws.Range("M60").Value = -19489
$ws.Range("N60").Value = -40784.4
$ws.Range("H61").Value = 29861.75
$ws.Range("J61").Value = 29861.75
$ws.Range("L61").Value = 29861.75
$ws.Range("N61").Value = -30557.75
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488
$ws.Range("H132").Value = 1964.3914
$ws.Range("I132").Value = 1599
$ws.Range("J132").Value = 2799.5715
$ws.Range("K132").Value = 4797
$ws.Range("L132").Value = 8398.7145
$ws.Range("M132").Value = -2267
$ws.Range("N132").Value = -13458.7145

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""
$ws.Range("H107").Value = 512.93335
$ws.Range("I107").Value = 585
$ws.Range("J107").Value = 368.8
$ws.Range("K107").Value = 1755
$ws.Range("L107").Value = 1106.4
$ws.Range("M107").Value = 165
$ws.Range("N107").Value = -4946.4
$ws.Range("H113").Value = 556187.8
$ws.Range("I113").Value = 1429153.2
$ws.Range("J113").Value = 664.36365
$ws.Range("K113").Value = 4287459.6
$ws.Range("L113").Value = 1993.09095
$ws.Range("M113").Value = -4285289.6
$ws.Range("N113").Value = -6333.09095
$ws.Range("H133").Value = 6333.75
$ws.Range("J133").Value = 6524.2856
$ws.Range("L133").Value = 19572.8568
$ws.Range("N133").Value = -29692.8568

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2250
$ws.Range("I132").Value = 1636.3334
$ws.Range("J132").Value = 3860.875
$ws.Range("K132").Value = 4909.0002
$ws.Range("L132").Value = 11582.625
$ws.Range("M132").Value = -2379.0002
$ws.Range("N132").Value = -16642.625

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 19333.334
$ws.Range("I21").Value = 19000
$ws.Range("K21").Value = 19000
$ws.Range("M21").Value = -18826
$ws.Range("H46").Value = 1590.5714
$ws.Range("I46").Value = 1421.4286
$ws.Range("J46").Value = 1928.8572
$ws.Range("K46").Value = 1421.4286
$ws.Range("L46").Value = 1928.8572
$ws.Range("M46").Value = -1233.4286
$ws.Range("N46").Value = -2304.8572
$ws.Range("H69").Value = 119862.2
$ws.Range("I69").Value = 9148
$ws.Range("J69").Value = 147540.75
$ws.Range("K69").Value = 9148
$ws.Range("L69").Value = 147540.75
$ws.Range("M69").Value = -8337
$ws.Range("N69").Value = -149162.75
$ws.Range("H72").Value = 119862.2
$ws.Range("I72").Value = 9148
$ws.Range("J72").Value = 147540.75
$ws.Range("K72").Value = 27444
$ws.Range("L72").Value = 442622.25
$ws.Range("M72").Value = -23388
$ws.Range("N72").Value = -450734.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 29006.5
$ws.Range("J25").Value = 29006.5
$ws.Range("L25").Value = 29006.5
$ws.Range("N25").Value = -29592.5
$ws.Range("H122").Value = 3397.0435
$ws.Range("I122").Value = 3183.0588
$ws.Range("J122").Value = 4003.3333
$ws.Range("K122").Value = 9549.1764
$ws.Range("L122").Value = 12009.9999
$ws.Range("M122").Value = -7099.1764
$ws.Range("N122").Value = -16909.9999
$ws.Range("H136").Value = 1183.5588
$ws.Range("I136").Value = 1174.7333
$ws.Range("K136").Value = 3524.199900000001
$ws.Range("M136").Value = -974.1999000000005
